$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 27
# Leading apostrophe forces the date-looking string to stay text, matching
# the other "Date" column cells in this sheet (inline string, not a date
# serial number).
$ws.Cells.Item($row, 1).Value = "'02/13/2026"
$ws.Cells.Item($row, 2).Value = 1553.712000000001
$ws.Cells.Item($row, 3).Value = 0.03185918625845714
$ws.Cells.Item($row, 4).Value = 50
